# Update the "想去人数" (want-to-go count) column F on the 展览 and 全部类型
# sheets to match the latest scraped totals.

$wb = $excel.ActiveWorkbook

# Row -> (old, new) value map for the "展览" sheet (F column)
$changesExhibition = @{
    4  = 323
    6  = 684
    7  = 271
    9  = 57
    12 = 3379
    13 = 107
    14 = 77
    18 = 575
    19 = 52
    20 = 680
    22 = 111
    25 = 59
    26 = 2425
    27 = 4959
    31 = 1272
    32 = 274
    33 = 2194
    35 = 484
    37 = 85
    38 = 155
    41 = 776
    42 = 25
    43 = 449
    45 = 456
}

# Row -> new value map for the "全部类型" sheet (F column) - row numbers are
# shifted by +1 from row 16 onward relative to the 展览 sheet because this
# sheet contains one additional entry.
$changesAllTypes = @{
    4  = 323
    6  = 684
    7  = 271
    9  = 57
    12 = 3379
    13 = 107
    14 = 77
    19 = 575
    20 = 52
    21 = 680
    23 = 111
    26 = 59
    27 = 2425
    28 = 4959
    32 = 1272
    33 = 274
    34 = 2194
    36 = 484
    38 = 85
    39 = 155
    42 = 776
    43 = 25
    44 = 449
    46 = 456
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $changesExhibition.Keys) {
    $wsExhibition.Range("F$row").Value = $changesExhibition[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $changesAllTypes.Keys) {
    $wsAllTypes.Range("F$row").Value = $changesAllTypes[$row]
}
